# Regenerate the handoff report: the localization job produced a new
# source-file GUID (60132e1a-266b-45f4-a93d-f0bd05f912eb) and new xliff
# fingerprints/timestamps. Update the workbook cells accordingly, keeping
# the existing hyperlink targets (they still point at the old commit/GUID
# in the repo - only the displayed file names and report data change).

$wb = $excel.ActiveWorkbook

$oldGuid = "64ae81cf-dbca-4c6e-93d5-684d8046acb4"
$newGuid = "60132e1a-266b-45f4-a93d-f0bd05f912eb"
$hyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/27372613fa499ad1d655d79f771064dd8bb2076e/e2e/$oldGuid.md"

# NOTE: this interpreter's PowerShell dialect does not bind named (-Param)
# arguments on user-defined functions, so Set-HyperlinkDisplay is invoked
# with positional arguments only.
function Set-HyperlinkDisplay {
    param($Worksheet, $CellRange, $DisplayText, $Address)

    $range = $Worksheet.Range($CellRange)
    $range.Hyperlinks.Delete()
    $Worksheet.Hyperlinks.Add($range, $Address, [Type]::Missing, [Type]::Missing, $DisplayText)
}

# ----- Overview sheet -----
$wsOverview = $wb.Worksheets.Item("Overview")

# A2: "<guid>.md" -> new guid
$wsOverview.Range("A2").Value = "$newGuid.md"

# B2: hyperlink display text "e2e\<guid>.md" -> new guid (keep same target URL)
Set-HyperlinkDisplay $wsOverview "B2" "e2e\$newGuid.md" $hyperlinkAddress

# G2: Latest HO Xliff Generate Date
$wsOverview.Range("G2").Value = "2016-08-22 05:06:33"

# ----- zh-cn sheet -----
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# A2: hyperlink display text "<guid>.md" -> new guid (keep same target URL)
Set-HyperlinkDisplay $wsZhCn "A2" "$newGuid.md" $hyperlinkAddress

# G2: Latest Handoff File
$wsZhCn.Range("G2").Value = "$newGuid.bf1b3947fcb901a189fffa81eb3cdc8557565fdb.zh-cn.xlf"

# H2: Latest Handoff Datetime
$wsZhCn.Range("H2").Value = "2016-08-22 05:06:28"

# ----- de-de sheet -----
$wsDeDe = $wb.Worksheets.Item("de-de")

# A2: hyperlink display text "<guid>.md" -> new guid (keep same target URL)
Set-HyperlinkDisplay $wsDeDe "A2" "$newGuid.md" $hyperlinkAddress

# G2: Latest Handoff File
$wsDeDe.Range("G2").Value = "$newGuid.bf1b3947fcb901a189fffa81eb3cdc8557565fdb.de-de.xlf"

# H2: Latest Handoff Datetime (mirrors Overview's Latest HO Xliff Generate Date)
$wsDeDe.Range("H2").Value = "2016-08-22 05:06:33"
